$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.502.21"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "1.825.96"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  -0.46%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "312.38"
$c.Style = "Normal"

$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("E6").Value = "  -0.26%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4255"
$c.Style = "Normal"

$ws.Range("E7").Value = "  -0.52%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3611"
$c.Style = "Normal"

$ws.Range("E8").Value = "  +0.81%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07206"
$c.Style = "Normal"

$ws.Range("E9").Value = "  -1.18%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.8620"
$c.Style = "Normal"

$ws.Range("E10").Value = "  -1.03%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "20.64"
$c.Style = "Normal"

$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "1.850.84"
$ws.Range("E12").Value = "  -3.01%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.388"
$c.Style = "Normal"

$ws.Range("E13").Value = "  +1.12%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.479"
$c.Style = "Normal"

$ws.Range("E14").Value = "  -1.02%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.06920"
$c.Style = "Normal"

$ws.Range("E15").Value = "  -1.15%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"

$ws.Range("E16").Value = "  -0.30%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "80.65"
$c.Style = "Normal"

$ws.Range("E17").Value = "  +1.29%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008929"
$c.Style = "Normal"

$ws.Range("E18").Value = "  -0.13%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"

$ws.Range("E19").Value = "  -0.39%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "15.35"
$c.Style = "Normal"

$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").Value = "27.378.74"
$ws.Range("E21").Value = "  -1.77%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.120"
$c.Style = "Normal"

$ws.Range("E22").Value = "  +2.51%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.88"
$c.Style = "Normal"

$ws.Range("E23").Value = "  +4.87%  "
$ws.Range("D24").Value = "2.031.26"
$ws.Range("E24").Value = "  -3.26%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.990"
$c.Style = "Normal"

$ws.Range("E25").Value = "  -0.05%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "155.05"
$c.Style = "Normal"

$ws.Range("E26").Value = "  -0.32%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "18.68"
$c.Style = "Normal"

$ws.Range("E27").Value = "  +1.04%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "5.152"
$c.Style = "Normal"

$ws.Range("E28").Value = "  -2.32%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "114.15"
$c.Style = "Normal"

$ws.Range("E29").Value = "  -5.47%  "
$ws.Range("E30").Value = "  -4.19%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.08873"
$c.Style = "Normal"

$ws.Range("E31").Value = "  -0.48%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.7506"
$c.Style = "Normal"

$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "2.974"
$c.Style = "Normal"

$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.537"
$c.Style = "Normal"

$ws.Range("E34").Value = "  +0.69%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.121"
$c.Style = "Normal"

$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("E36").Value = "  -0.27%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.086"
$c.Style = "Normal"

$ws.Range("E37").Value = "  -1.20%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.05282"
$c.Style = "Normal"

$ws.Range("E38").Value = "  -2.65%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01921"
$c.Style = "Normal"

$ws.Range("E39").Value = "  -0.40%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.796"
$c.Style = "Normal"

$ws.Range("E40").Value = "  -1.10%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.5072"
$c.Style = "Normal"

$ws.Range("E41").Value = "  +0.07%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1656"
$c.Style = "Normal"

$ws.Range("E42").Value = "  -0.57%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "6.452"
$c.Style = "Normal"

$ws.Range("E43").Value = "  -2.64%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.341"
$c.Style = "Normal"

$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("E45").Value = "  +1.00%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "106.28"
$c.Style = "Normal"

$ws.Range("E46").Value = "  -0.13%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.06465"
$c.Style = "Normal"

$ws.Range("E47").Value = "  -1.21%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.4675"
$c.Style = "Normal"

$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("E50").Value = "  -0.51%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "63.87"
$c.Style = "Normal"

$ws.Range("E51").Value = "  -0.75%  "